$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: re-order existing rows (3-73) whose betting-odds cycle (columns
# F:V = home..url_partida) was re-sorted in the upstream scrape/export.
# Columns A:E (Indice, pais, torneio, temporada, data_partida) are not
# touched because they are fixed to the row position.
# ---------------------------------------------------------------------------

$old3 = $ws.Range("F3:V3").Value2
$old4 = $ws.Range("F4:V4").Value2
$old5 = $ws.Range("F5:V5").Value2
$ws.Range("F3:V3").Value2 = $old4
$ws.Range("F4:V4").Value2 = $old5
$ws.Range("F5:V5").Value2 = $old3

$old11 = $ws.Range("F11:V11").Value2
$old12 = $ws.Range("F12:V12").Value2
$ws.Range("F11:V11").Value2 = $old12
$ws.Range("F12:V12").Value2 = $old11

$old17 = $ws.Range("F17:V17").Value2
$old19 = $ws.Range("F19:V19").Value2
$old18 = $ws.Range("F18:V18").Value2
$ws.Range("F17:V17").Value2 = $old19
$ws.Range("F19:V19").Value2 = $old18
$ws.Range("F18:V18").Value2 = $old17

$old29 = $ws.Range("F29:V29").Value2
$old31 = $ws.Range("F31:V31").Value2
$ws.Range("F29:V29").Value2 = $old31
$ws.Range("F31:V31").Value2 = $old29

$old35 = $ws.Range("F35:V35").Value2
$old37 = $ws.Range("F37:V37").Value2
$ws.Range("F35:V35").Value2 = $old37
$ws.Range("F37:V37").Value2 = $old35

$old38 = $ws.Range("F38:V38").Value2
$old39 = $ws.Range("F39:V39").Value2
$ws.Range("F38:V38").Value2 = $old39
$ws.Range("F39:V39").Value2 = $old38

$old44 = $ws.Range("F44:V44").Value2
$old46 = $ws.Range("F46:V46").Value2
$old45 = $ws.Range("F45:V45").Value2
$ws.Range("F44:V44").Value2 = $old46
$ws.Range("F46:V46").Value2 = $old45
$ws.Range("F45:V45").Value2 = $old44

$old47 = $ws.Range("F47:V47").Value2
$old48 = $ws.Range("F48:V48").Value2
$ws.Range("F47:V47").Value2 = $old48
$ws.Range("F48:V48").Value2 = $old47

$old50 = $ws.Range("F50:V50").Value2
$old51 = $ws.Range("F51:V51").Value2
$old52 = $ws.Range("F52:V52").Value2
$ws.Range("F50:V50").Value2 = $old51
$ws.Range("F51:V51").Value2 = $old52
$ws.Range("F52:V52").Value2 = $old50

$old56 = $ws.Range("F56:V56").Value2
$old58 = $ws.Range("F58:V58").Value2
$ws.Range("F56:V56").Value2 = $old58
$ws.Range("F58:V58").Value2 = $old56

$old62 = $ws.Range("F62:V62").Value2
$old64 = $ws.Range("F64:V64").Value2
$ws.Range("F62:V62").Value2 = $old64
$ws.Range("F64:V64").Value2 = $old62

$old68 = $ws.Range("F68:V68").Value2
$old69 = $ws.Range("F69:V69").Value2
$ws.Range("F68:V68").Value2 = $old69
$ws.Range("F69:V69").Value2 = $old68

# ---------------------------------------------------------------------------
# Part 2: append 3 new scraped matches (rows 74-76) for the matchday of
# 28/10/2023-29/10/2023, pushing the used range from A1:V73 to A1:V76.
# ---------------------------------------------------------------------------

# Clone formatting (styles) of the last existing data row onto the new rows
# so the new cells share the same number formats / borders as the rest of
# the table (bold+border style for column A, date style for column E).
$ws.Range("A73:V73").Copy() | Out-Null
$ws.Range("A74:V76").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(74,1).Value2 = 73
$ws.Cells.Item(74,2).Value2 = "austria"
$ws.Cells.Item(74,3).Value2 = "bundesliga"
$ws.Cells.Item(74,4).Value2 = "2023-2024"
$ws.Cells.Item(74,5).Value2 = 45234.70833333334
$ws.Cells.Item(74,6).Value2 = "Tirol"
$ws.Cells.Item(74,7).Value2 = 0
$ws.Cells.Item(74,8).Value2 = "Salzburg"
$ws.Cells.Item(74,9).Value2 = 2
$ws.Cells.Item(74,10).Value2 = 8.32
$ws.Cells.Item(74,11).Value2 = "28/10/2023 18:12"
$ws.Cells.Item(74,12).Value2 = 11.61
$ws.Cells.Item(74,13).Value2 = "04/11/2023 16:44"
$ws.Cells.Item(74,14).Value2 = 5.55
$ws.Cells.Item(74,15).Value2 = "28/10/2023 18:12"
$ws.Cells.Item(74,16).Value2 = 6.58
$ws.Cells.Item(74,17).Value2 = "04/11/2023 16:48"
$ws.Cells.Item(74,18).Value2 = 1.36
$ws.Cells.Item(74,19).Value2 = "28/10/2023 18:12"
$ws.Cells.Item(74,20).Value2 = 1.25
$ws.Cells.Item(74,21).Value2 = "04/11/2023 14:52"
$ws.Cells.Item(74,22).Value2 = "https://www.betexplorer.com/football/austria/bundesliga/tirol-salzburg/CjhR9jDC/"

$ws.Cells.Item(75,1).Value2 = 74
$ws.Cells.Item(75,2).Value2 = "austria"
$ws.Cells.Item(75,3).Value2 = "bundesliga"
$ws.Cells.Item(75,4).Value2 = "2023-2024"
$ws.Cells.Item(75,5).Value2 = 45234.70833333334
$ws.Cells.Item(75,6).Value2 = "Austria Vienna"
$ws.Cells.Item(75,7).Value2 = 1
$ws.Cells.Item(75,8).Value2 = "A. Lustenau"
$ws.Cells.Item(75,9).Value2 = 0
$ws.Cells.Item(75,10).Value2 = 1.37
$ws.Cells.Item(75,11).Value2 = "29/10/2023 14:42"
$ws.Cells.Item(75,12).Value2 = 1.28
$ws.Cells.Item(75,13).Value2 = "04/11/2023 16:54"
$ws.Cells.Item(75,14).Value2 = 5.29
$ws.Cells.Item(75,15).Value2 = "29/10/2023 14:42"
$ws.Cells.Item(75,16).Value2 = 6.23
$ws.Cells.Item(75,17).Value2 = "04/11/2023 16:56"
$ws.Cells.Item(75,18).Value2 = 7.99
$ws.Cells.Item(75,19).Value2 = "29/10/2023 14:42"
$ws.Cells.Item(75,20).Value2 = 10.31
$ws.Cells.Item(75,21).Value2 = "04/11/2023 16:56"
$ws.Cells.Item(75,22).Value2 = "https://www.betexplorer.com/football/austria/bundesliga/austria-vienna-a-lustenau/CGNMChTg/"

$ws.Cells.Item(76,1).Value2 = 75
$ws.Cells.Item(76,2).Value2 = "austria"
$ws.Cells.Item(76,3).Value2 = "bundesliga"
$ws.Cells.Item(76,4).Value2 = "2023-2024"
$ws.Cells.Item(76,5).Value2 = 45234.70833333334
$ws.Cells.Item(76,6).Value2 = "Wolfsberger AC"
$ws.Cells.Item(76,7).Value2 = 4
$ws.Cells.Item(76,8).Value2 = "A. Klagenfurt"
$ws.Cells.Item(76,9).Value2 = 0
$ws.Cells.Item(76,10).Value2 = 2.25
$ws.Cells.Item(76,11).Value2 = "28/10/2023 18:12"
$ws.Cells.Item(76,12).Value2 = 2.09
$ws.Cells.Item(76,13).Value2 = "04/11/2023 16:58"
$ws.Cells.Item(76,14).Value2 = 3.63
$ws.Cells.Item(76,15).Value2 = "28/10/2023 18:12"
$ws.Cells.Item(76,16).Value2 = 3.68
$ws.Cells.Item(76,17).Value2 = "04/11/2023 16:58"
$ws.Cells.Item(76,18).Value2 = 3.22
$ws.Cells.Item(76,19).Value2 = "28/10/2023 18:12"
$ws.Cells.Item(76,20).Value2 = 3.59
$ws.Cells.Item(76,21).Value2 = "04/11/2023 16:58"
$ws.Cells.Item(76,22).Value2 = "https://www.betexplorer.com/football/austria/bundesliga/wolfsberger-ac-a-klagenfurt/bDJIDYDm/"
